$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 12551
    3 = 12551
    4 = 12085
    5 = 11929
    6 = 11418
    7 = 11164
    8 = 11164
    9 = 11164
    10 = 10850
    11 = 10850
    12 = 10850
    13 = 10273
    14 = 10273
    15 = 10273
    16 = 10273
    17 = 10273
    18 = 10273
    19 = 10273
    20 = 10273
    21 = 10273
    22 = 10273
    23 = 8597
    24 = 8597
    25 = 8072
    26 = 8072
    27 = 8072
    28 = 8072
    29 = 8072
    30 = 8072
    31 = 8072
    32 = 8072
    33 = 8072
    34 = 7760
    35 = 7760
    36 = 7760
    37 = 7760
    38 = 7682
    39 = 7682
    40 = 7682
    41 = 7682
    42 = 7682
    43 = 7682
    44 = 7682
    45 = 7682
    46 = 7682
    47 = 7682
    48 = 7682
    49 = 7682
    50 = 7682
    51 = 7682
    52 = 7682
    53 = 7682
    54 = 7657
    55 = 7657
    56 = 7657
    57 = 7657
    58 = 7657
    59 = 7657
    60 = 7657
    61 = 7657
    62 = 7657
    63 = 7657
    64 = 7657
    65 = 7657
    66 = 7657
    67 = 7657
    68 = 7657
    69 = 7657
    70 = 7657
    71 = 7657
    72 = 7657
    73 = 7657
    74 = 7657
    75 = 7657
    76 = 7657
    77 = 7657
    78 = 7657
    79 = 7657
    80 = 7657
    81 = 7657
    82 = 7657
    83 = 7657
    84 = 7657
    85 = 7657
    86 = 7657
    87 = 7657
    88 = 7657
    89 = 7657
    90 = 7657
    91 = 7657
    92 = 7657
    93 = 7657
    94 = 7657
    95 = 7657
    96 = 7573
    97 = 7573
    98 = 7573
    99 = 7573
    100 = 7573
    101 = 7573
    102 = 7573
    103 = 7573
    104 = 7573
    105 = 7573
    106 = 7573
    107 = 7573
    108 = 7573
    109 = 7573
    110 = 7573
    111 = 7573
    112 = 7573
    113 = 7573
    114 = 7573
    115 = 7573
    116 = 7573
    117 = 7573
    118 = 7573
    119 = 7573
    120 = 7573
    121 = 7573
    122 = 7573
    123 = 7573
    124 = 7573
    125 = 7573
    126 = 7573
    127 = 7573
    128 = 7573
    129 = 7573
    130 = 7573
    131 = 7573
    132 = 7573
    133 = 7573
    134 = 7573
    135 = 7573
    136 = 7573
    137 = 7573
    138 = 7573
    139 = 7573
    140 = 7573
    141 = 7573
    142 = 7573
    143 = 7573
    144 = 7573
    145 = 7573
    146 = 7573
    147 = 7573
    148 = 7573
    149 = 7573
    150 = 7573
    151 = 7573
    152 = 7573
    153 = 7573
    154 = 7573
    155 = 7573
    156 = 7573
    157 = 7573
    158 = 7573
    159 = 7573
    160 = 7573
    161 = 7573
    162 = 7573
    163 = 7573
    164 = 7573
    165 = 7573
    166 = 7573
    167 = 7573
    168 = 7573
    169 = 7573
    170 = 7573
    171 = 7573
    172 = 7573
    173 = 7573
    174 = 7573
    175 = 7573
    176 = 7573
    177 = 7573
    178 = 7573
    179 = 7573
    180 = 7573
    181 = 7573
    182 = 7573
    183 = 7573
    184 = 7573
    185 = 7573
    186 = 7573
    187 = 7573
    188 = 7573
    189 = 7573
    190 = 7573
    191 = 7573
    192 = 7573
    193 = 7573
    194 = 7573
    195 = 7573
    196 = 7573
    197 = 7573
    198 = 7573
    199 = 7573
    200 = 7573
    201 = 7573
    202 = 7573
    203 = 7573
    204 = 7573
    205 = 7573
    206 = 7573
    207 = 7573
    208 = 7573
    209 = 7573
    210 = 7573
    211 = 7573
    212 = 7573
    213 = 7573
    214 = 7573
    215 = 7573
    216 = 7573
    217 = 7573
    218 = 7573
    219 = 7573
    220 = 7573
    221 = 7573
    222 = 7573
    223 = 7573
    224 = 7573
    225 = 7573
    226 = 7573
    227 = 7573
    228 = 7573
    229 = 7573
    230 = 7573
    231 = 7573
    232 = 7573
    233 = 7573
    234 = 7573
    235 = 7573
    236 = 7573
    237 = 7573
    238 = 7573
    239 = 7573
    240 = 7573
    241 = 7573
    242 = 7573
    243 = 7573
    244 = 7573
    245 = 7573
    246 = 7573
    247 = 7573
    248 = 7573
    249 = 7573
    250 = 7573
    251 = 7573
    252 = 7573
}

foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 3).Value = $values[$r]
}

Write-Output "Updated $($values.Count) cells in column C"